# "Generate Report for Handback"
#
# The handback CI step discovered that the zh-cn and de-de targets came
# back in sync with en-US, so it:
#   - flips the shared "Status" text from "Ready for handoff" to
#     "Handed back: in sync with en-US" (this text is shown on the
#     Overview rollup sheet as well as on each language sheet)
#   - records the generated target (.md) and handback (.xlf) file names
#   - stamps the handback datetime (distinct per language)
#   - links the new target-file cell back to the source doc, just like
#     the existing "Source File Name" link

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eca2003561992ffb0f1a5713d52595e14a265753/e2e/b742e51e-0df9-44be-a16e-1a022713b4da.md"
$mdName = "b742e51e-0df9-44be-a16e-1a022713b4da.md"
$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

function Set-HandbackRow {
    param($ws, [string]$handbackFile, [string]$handbackDateTime)

    # Status (column C)
    $ws.Range("C2").Value = $statusText

    # Latest Target File (column I) - gets the handed-back doc name, linked
    # back to the source markdown file, same as column A's hyperlink.
    $ws.Range("I2").Value = $mdName
    $ws.Hyperlinks.Add($ws.Range("I2"), $mdUrl, "", "", $mdName) | Out-Null

    # Latest Handback File (column J)
    $ws.Range("J2").Value = $handbackFile

    # Latest Handback DateTime (column K)
    $ws.Range("K2").Value = $handbackDateTime

    # Widen the columns so the longer status text and the long file names
    # are not truncated - matches the other filename columns (A, G) which
    # are already fixed at a width of 40.
    $ws.Columns.Item(3).ColumnWidth = 29.14
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

Set-HandbackRow $wsZhCn `
    "b742e51e-0df9-44be-a16e-1a022713b4da.25b5561527b3d21c1c8e4884f0d0b37954a58ac1.zh-cn.xlf" `
    "2016-08-15 20:55:01"

Set-HandbackRow $wsDeDe `
    "b742e51e-0df9-44be-a16e-1a022713b4da.25b5561527b3d21c1c8e4884f0d0b37954a58ac1.de-de.xlf" `
    "2016-08-15 20:55:15"

# The Overview sheet mirrors the same Status text for each language.
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14
